# Fruta / hortaliza, semanal
# Insert 2 new daily price rows for "Naranja" (Fukumoto variety) at row 670,
# pushing the existing rows 670-748 down to 672-750.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 670 (existing row 670 and everything
# below shifts down by 2 rows).
$ws.Rows.Item(670).Resize(2).Insert()

# --- New row 670 ---------------------------------------------------------
$ws.Cells.Item(670, 1).Value = 7
$ws.Cells.Item(670, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(670, 3).Value = "Ñuble"
$ws.Cells.Item(670, 4).Value = 45142
$ws.Cells.Item(670, 5).Value = 16
$ws.Cells.Item(670, 6).Value = "Fruta"
$ws.Cells.Item(670, 7).Value = 100102
$ws.Cells.Item(670, 8).Value = "Cítricos"
$ws.Cells.Item(670, 9).Value = 100102005
$ws.Cells.Item(670, 10).Value = "Naranja"
$ws.Cells.Item(670, 11).Value = "Fukumoto"
$ws.Cells.Item(670, 12).Value = "Primera"
$ws.Cells.Item(670, 13).Value = 80
$ws.Cells.Item(670, 14).Value = 8000
$ws.Cells.Item(670, 15).Value = 8000
$ws.Cells.Item(670, 16).Value = 8000
$ws.Cells.Item(670, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(670, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(670, 19).Value = 533
$ws.Cells.Item(670, 20).Value = 15

# --- New row 671 ---------------------------------------------------------
$ws.Cells.Item(671, 1).Value = 7
$ws.Cells.Item(671, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(671, 3).Value = "Ñuble"
$ws.Cells.Item(671, 4).Value = 45142
$ws.Cells.Item(671, 5).Value = 16
$ws.Cells.Item(671, 6).Value = "Fruta"
$ws.Cells.Item(671, 7).Value = 100102
$ws.Cells.Item(671, 8).Value = "Cítricos"
$ws.Cells.Item(671, 9).Value = 100102005
$ws.Cells.Item(671, 10).Value = "Naranja"
$ws.Cells.Item(671, 11).Value = "Fukumoto"
$ws.Cells.Item(671, 12).Value = "Segunda"
$ws.Cells.Item(671, 13).Value = 80
$ws.Cells.Item(671, 14).Value = 6000
$ws.Cells.Item(671, 15).Value = 6000
$ws.Cells.Item(671, 16).Value = 6000
$ws.Cells.Item(671, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(671, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(671, 19).Value = 400
$ws.Cells.Item(671, 20).Value = 15
